$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '63.929.75'
Set-TextValue 'E2' '  +0.07%  '
Set-TextValue 'D3' '3.102.85'
Set-TextValue 'E3' '  +1.12%  '
Set-TextValue 'E4' '  +0.20%  '
Set-TextValue 'D5' '558.19'
Set-TextValue 'E5' '  +1.46%  '
Set-TextValue 'D6' '138.67'
Set-TextValue 'E6' '  +0.41%  '
Set-TextValue 'E7' '  +0.39%  '
Set-TextValue 'D8' '3.097.24'
Set-TextValue 'E8' '  +1.23%  '
Set-TextValue 'D9' '0.491'
Set-TextValue 'E9' '  +0.83%  '
Set-TextValue 'E10' '  +5.03%  '
Set-TextValue 'E11' '  -0.21%  '
Set-TextValue 'E12' '  +0.91%  '
Set-TextValue 'E13' '  +0.95%  '
Set-TextValue 'E14' '  +0.30%  '
Set-TextValue 'D15' '3.614.77'
Set-TextValue 'E15' '  +2.34%  '
Set-TextValue 'D16' '63.995.70'
Set-TextValue 'E16' '  +0.16%  '
Set-TextValue 'E17' '  +0.44%  '
Set-TextValue 'D18' '3.112.75'
Set-TextValue 'E18' '  +2.84%  '
Set-TextValue 'D19' '505.56'
Set-TextValue 'E19' '  +4.51%  '
Set-TextValue 'E20' '  +1.60%  '
Set-TextValue 'E21' '  +1.29%  '
Set-TextValue 'E22' '  +3.33%  '
Set-TextValue 'E23' '  +1.11%  '
Set-TextValue 'D24' '12.43'
Set-TextValue 'E24' '  +0.78%  '
Set-TextValue 'D25' '77.90'
Set-TextValue 'E25' '  +0.27%  '
Set-TextValue 'D26' '1.00'
Set-TextValue 'E26' '  -0.26%  '
Set-TextValue 'D27' '2.78'
Set-TextValue 'E27' '  +3.71%  '
Set-TextValue 'D28' '8.46'
Set-TextValue 'E28' '  +6.40%  '
Set-TextValue 'E29' '  -0.33%  '
Set-TextValue 'D30' '0.999'
Set-TextValue 'E30' '  +0.35%  '
Set-TextValue 'D31' '26.24'
Set-TextValue 'E31' '  +1.93%  '
Set-TextValue 'E32' '  -2.49%  '
Set-TextValue 'E33' '  +0.35%  '
Set-TextValue 'D34' '543.23'
Set-TextValue 'E34' '  -6.92%  '
Set-TextValue 'D35' '55.10'
Set-TextValue 'E35' '  +6.16%  '
Set-TextValue 'D36' '5.91'
Set-TextValue 'E36' '  -0.58%  '
Set-TextValue 'D37' '5.22'
Set-TextValue 'E37' '  -2.93%  '
Set-TextValue 'D38' '0.0416'
Set-TextValue 'E38' '  +4.39%  '
Set-TextValue 'D39' '0.0801'
Set-TextValue 'E39' '  +1.96%  '
Set-TextValue 'D40' '3.069.69'
Set-TextValue 'E40' '  +4.50%  '
Set-TextValue 'E41' '  +1.22%  '
Set-TextValue 'E42' '  -0.47%  '
Set-TextValue 'D43' '2.63'
Set-TextValue 'E43' '  -9.84%  '
Set-TextValue 'D44' '0.256'
Set-TextValue 'E44' '  +5.92%  '
Set-TextValue 'E45' '  +0.04%  '
Set-TextValue 'D46' '2.11'
Set-TextValue 'E46' '  +1.04%  '
Set-TextValue 'D47' '120.88'
Set-TextValue 'E47' '  +1.49%  '
Set-TextValue 'D48' '24.33'
Set-TextValue 'E48' '  -2.14%  '
Set-TextValue 'E49' '  -0.26%  '
Set-TextValue 'E50' '  -4.05%  '
Set-TextValue 'E51' '  -0.80%  '
